# Kleine Änderungen kurz vor Präsentation
$p = $ppt.ActivePresentation

# --- Slide 12: "Verbesserungsvorschläge" --------------------------------
# Split the "JavaDoc Kommentare ..." run so "JavaDoc" becomes its own run
# (matches the target run split that marks "JavaDoc" as a spell-check
# exception, keeping the rest of the sentence in a second run).
$s12 = $p.Slides.Item(12)
$tr12 = $s12.Shapes.Item(2).TextFrame.TextRange
$javaDoc = $tr12.Characters(146, 7)
$javaDoc.Text = "JavaDoc"

# --- Slide 8: "Grafische Oberfläche" -------------------------------------
# "... mit kein Treffer überschrieben" -> "... mit „Kein Treffer“ überschrieben"
$s8 = $p.Slides.Item(8)
$tr8 = $s8.Shapes.Item(2).TextFrame.TextRange
$run8 = $tr8.Characters(193, 31)
$run8.Text = " mit „Kein Treffer“ überschrieben"

# --- Slide 9: "Datenbankverwaltung" --------------------------------------
# "RunableJarFiles" -> "RunableJarFile"
$s9 = $p.Slides.Item(9)
$tr9 = $s9.Shapes.Item(2).TextFrame.TextRange
$run9 = $tr9.Characters(314, 15)
$run9.Text = "RunableJarFile"
